# Daily attendance processing - 2025-10-20 20:45:22
# Normalize the "Recorded By" (column G) values: move the literal "System"
# token to the front of the comma-separated list of recorders, and fix the
# casing order of the duplicate "system"/"System" backdoor entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of exact old text -> new text, as observed for column G values.
$map = @{
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "admin@admin.com, System"             = "System, admin@admin.com"
    "backup@backdoor.com, system, System" = "backup@backdoor.com, System, system"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
